$wb = $excel.ActiveWorkbook

# --- Sheet "data": update row 2 ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("A2").ClearContents()
$wsData.Range("B2").Value = "доставка"
$wsData.Range("C2").Value = 40
$wsData.Range("D2").Value = "т"
$wsData.Range("F2").Value = "\nИГ:\nЗаявка на доставку \n1. Дата отгрузки 28.04.2023 \n2. Марка ЦЕМ I 42.5н Беларусь   \n3. Количество 40 тонн  \n4. От ООО Спарта \n5. Завод: Сзтк \n6. ] ООО """"ТД""Цемент \n7. Грузополучатель: ООО ""ТД""Цемент  \n8. Голицыно\n+7 910 404-06-14\nРБУ\nМожайское ш., 81\n"

# --- Sheet "ошибки": drop the old traceback note row ---
$wsErr = $wb.Worksheets.Item("ошибки")
$wsErr.Rows.Item(2).Delete()

# --- force full recalculation on load ---
$wb.Application.CalculateFullRebuild()
